$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has two rows:
#   A1 = 0            (bold, bordered/centered style)
#   A2 = <questions text, shared string>
# The target keeps only the questions text (reformatted as pretty JSON)
# in A1, with the default (unstyled) formatting, and removes the old A1
# row entirely (row shift up).

# Capture the (new, reformatted) text that should end up in A1.
$questionsText = @'
questions = [
    {
        "title": "Describe yourself as you generally are now, not as you wish to be in the future. Describe yourself as you honestly see yourself, in relation to your peers.",
        "ques_type": 2,
        "options": [
            "I rarely display signs of affection.",
            "I am not afraid to display signs of affection."
        ],
        "score": "I rarely display signs of affection."
    },
    {
        "title": "Describe yourself as you generally are now, not as you wish to be in the future. Describe yourself as you honestly see yourself, in relation to your peers.",
        "ques_type": 2,
        "options": [
            "Friendship over fairness.",
            "Fairness over friendship."
        ],
        "score": "Friendship over fairness."
    },
    {
        "title": "Describe yourself as you generally are now, not as you wish to be in the future. Describe yourself as you honestly see yourself, in relation to your peers.",
        "ques_type": 2,
        "options": [
            "I speak my mind about other people\u2019s lives.",
            "I am withdrawn and somewhat ambiguous in my communication."
        ],
        "score": "I speak my mind about other people\u2019s lives."
    },
    {
        "title": "Describe yourself as you generally are now, not as you wish to be in the future. Describe yourself as you honestly see yourself, in relation to your peers.",
        "ques_type": 2,
        "options": [
            "I am determined in times of hardship.",
            "I withdraw and am demotivated when faced with hardship."
        ],
        "score": "I am determined in times of hardship."
    }
]
'@

# Write the freshly reformatted text (pretty-printed JSON instead of the
# old Python-literal style single line) into the cell that currently
# holds the questions text (A2).
$ws.Range("A2").Value = $questionsText

# Setting a value with embedded line breaks causes Excel to auto-fit the
# row height (adding an explicit custom height). Restore the row to the
# sheet's standard height so no stray height formatting is left behind.
$ws.Rows.Item(2).AutoFit()

# Delete the first row (value 0 / bold+bordered style). This shifts the
# questions-text row up into A1, so only the questions text remains, as
# a single row/cell, with the default (unstyled) formatting.
$ws.Rows.Item(1).Delete()
